$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell reference -> new value, applied as literal text (matches the
# source data which stores prices/volumes as inline strings, e.g.
# "1.000" or "29.145.16", not as numbers).
$updates = [ordered]@{
    'D2' = '29.145.16'
    'E2' = '  -3.28%  '
    'D3' = '1.849.81'
    'E3' = '  -2.24%  '
    'D4' = '1.000'
    'E4' = '  -0.04%  '
    'D5' = '0.7021'
    'E5' = '  -5.11%  '
    'D6' = '238.18'
    'D7' = '1.000'
    'E7' = '  -0.01%  '
    'D8' = '0.3035'
    'E8' = '  -4.35%  '
    'D9' = '0.07462'
    'E9' = '  +3.53%  '
    'D10' = '23.32'
    'E10' = '  -6.47%  '
    'D11' = '0.08129'
    'E11' = '  -2.63%  '
    'D12' = '1.843.29'
    'E12' = '  -5.42%  '
    'D13' = '0.7243'
    'E13' = '  -4.76%  '
    'D14' = '5.214'
    'E14' = '  -4.36%  '
    'D15' = '88.65'
    'E15' = '  -4.91%  '
    'D16' = '29.136.42'
    'E16' = '  -3.40%  '
    'D17' = '5.757'
    'E17' = '  -6.82%  '
    'D18' = '236.51'
    'E18' = '  -5.72%  '
    'D19' = '13.05'
    'E19' = '  -4.31%  '
    'E20' = '  -3.41%  '
    'D21' = '0.9999'
    'E21' = '  -0.08%  '
    'D22' = '2.094.04'
    'E22' = '  -4.33%  '
    'D23' = '1.000'
    'E23' = '  -0.02%  '
    'D24' = '7.540'
    'E24' = '  -5.24%  '
    'D25' = '8.984'
    'E25' = '  -3.61%  '
    'D26' = '160.85'
    'E26' = '  -2.31%  '
    'D27' = '0.1449'
    'E27' = '  -8.49%  '
    'D28' = '18.02'
    'E28' = '  -4.09%  '
    'D29' = '1.962'
    'E29' = '  -5.10%  '
    'D30' = '1.398'
    'E30' = '  -5.99%  '
    'D31' = '4.515'
    'E31' = '  -1.68%  '
    'E32' = '  -3.18%  '
    'D33' = '3.959'
    'E33' = '  -6.05%  '
    'D34' = '0.05142'
    'E34' = '  -4.33%  '
    'D35' = '1.183'
    'E35' = '  -5.92%  '
    'D36' = '1.033'
    'E36' = '  +2.56%  '
    'D37' = '0.6976'
    'E37' = '  -10.20%  '
    'D38' = '2.659'
    'E38' = '  -2.53%  '
    'D39' = '0.01862'
    'E39' = '  -5.24%  '
    'D40' = '2.676'
    'E40' = '  -3.22%  '
    'D41' = '0.9384'
    'E41' = '  +6.25%  '
    'B42' = 'FraxShare'
    'C42' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D42' = '5.985'
    'E42' = '  -1.77%  '
    'B43' = 'Maker'
    'C43' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D43' = '1.078.68'
    'E43' = '  -1.88%  '
    'D44' = '0.4274'
    'E44' = '  -6.51%  '
    'D45' = '69.70'
    'E45' = '  -4.11%  '
    'E46' = '  -0.17%  '
    'D47' = '101.95'
    'E47' = '  -2.56%  '
    'D48' = '1.737'
    'E48' = '  -7.01%  '
    'D49' = '1.985.26'
    'E49' = '  -4.31%  '
    'B50' = 'Aptos'
    'C50' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D50' = '7.018'
    'E50' = '  -7.70%  '
    'B51' = 'EnergySwap'
    'C51' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D51' = '9.110'
    'E51' = '  -5.45%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force a text-literal write so numeric-looking strings (e.g. "1.000",
    # "0.9999") are not auto-coerced into numbers by value assignment,
    # then clear the temporary format so the cell ends up styled exactly
    # as it started (no NumberFormat residue on the cell).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.ClearFormats()
}

Write-Host ("Updated " + $updates.Count + " cells")
